$d = $word.ActiveDocument

$tbl = $d.Tables.Item(3)
$row = $tbl.Rows.Item($tbl.Rows.Count)
$cell = $row.Cells.Item($row.Cells.Count)
$rng = $cell.Range
Write-Output "cellrange:"
Write-Output $rng.Start
Write-Output $rng.End
Write-Output ("[" + $rng.Text + "]")

$para = $cell.Range.Paragraphs.Item(1)
$prng = $para.Range
Write-Output "pararange:"
Write-Output $prng.Start
Write-Output $prng.End
Write-Output ("[" + $prng.Text + "]")

$prng.InsertAfter("Costante Marco")
Write-Output "after insert:"
Write-Output $prng.Start
Write-Output $prng.End
Write-Output ("[" + $prng.Text + "]")

$prng.NoProofing = 1
Write-Output "done"
